$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(46051, 11974.5398612447, 11586.069059587, 17595.86, 7583.9167304747, 65.5885745859029),
    @(46052, 12219.3067193567, 11133.8787734538, 12075.86, 8239.66844358911, 304.07030071012),
    @(46053, 4963.57835265955, 7981.63782676138, 12075.86, 8108.69528621942, 167.269713040867),
    @(46054, 5376.34744431836, 8018.07757703374, 9791.86, 8112.38400172761, 264.108399115056),
    @(46055, 11845.3920979544, 11300.8279437159, 9791.86, 7924.5943990593, 393.065097615635),
    @(46056, 11503.27393928, 10652.7300015177, 9791.86, 7642.42654437646, 354.304022745588),
    @(46057, 11503.27393928, 10399.6452063125, 9791.86, 7642.42654437646, 343.758822945372),
    @(46058, 11503.27393928, 10655.9753097572, 9791.86, 7642.42654437646, 354.439243922237),
    @(46059, 11503.27393928, 9869.35279423597, 9791.86, 7642.42654437646, 321.663305775518),
    @(46060, 4875.64205484986, 7066.0382017956, 9791.86, 7273.63627175629, 189.492269731329),
    @(46061, 4779.1830818516, 7092.2104883314, 9791.86, 7266.27717281881, 190.276152547926),
    @(46062, 11115.7963561476, 10506.0722017952, 9791.86, 7172.64659512332, 328.61911653827),
    @(46063, 11115.7963561476, 10732.083983898, 9791.86, 7172.64659512332, 338.036274125888),
    @(46064, 11115.7963561476, 10774.6371550902, 9791.86, 7172.64659512332, 339.809322925565)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
